$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated total_risk (R) and total_risk_resp (S) values per newest airtoxics NATA data
$ws.Range("R2").Value = 37.3134328358209
$ws.Range("S2").Value = 0.440298507462687

$ws.Range("R3").Value = 52.972972972973
$ws.Range("S3").Value = 0.510810810810811

$ws.Range("R4").Value = 31.304347826087
$ws.Range("S4").Value = 0.411159420289855

$ws.Range("R5").Value = 29.7540983606557
$ws.Range("S5").Value = 0.334426229508197

$ws.Range("R6").Value = 20.8739495798319
$ws.Range("S6").Value = 0.295294117647059

$ws.Range("R7").Value = 20
$ws.Range("S7").Value = 0.275

$ws.Range("S8").Value = 0.2

$ws.Range("R9").Value = 54.0677966101695
$ws.Range("S9").Value = 0.413559322033898

$ws.Range("R10").Value = 53.9130434782609
$ws.Range("S10").Value = 0.427536231884058

$ws.Range("R11").Value = 34.4117647058824
$ws.Range("S11").Value = 0.408823529411765

$ws.Range("R12").Value = 73.3333333333333
$ws.Range("S12").Value = 0.311111111111111

$ws.Range("R13").Value = 18.7878787878788
$ws.Range("S13").Value = 0.181818181818182

$ws.Range("R14").Value = 30.1754385964912
$ws.Range("S14").Value = 0.331578947368421

$ws.Range("R15").Value = 30.4285714285714
$ws.Range("S15").Value = 0.345714285714286
